$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 518 (La Araucania / Vega Modelo
# de Temuco / Berenjena data set). Insert a whole row there so every
# subsequent record shifts down by one (old row 518 -> new row 519, ...,
# old row 544 -> new row 545).
$ws.Rows.Item(518).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A518").Value = 10
$ws.Range("B518").Value = "Vega Modelo de Temuco"
$ws.Range("C518").Value = "La Araucanía"
$ws.Range("D518").Value = 45267
$ws.Range("E518").Value = 9
$ws.Range("F518").Value = 100112001
$ws.Range("G518").Value = "Berenjena"
$ws.Range("H518").Value = "Sin especificar"
$ws.Range("I518").Value = "Primera"
$ws.Range("J518").Value = 210
$ws.Range("K518").Value = 14000
$ws.Range("L518").Value = 15000
$ws.Range("M518").Value = 14429
$ws.Range("N518").Value = '$/caja 40 unidades'
$ws.Range("O518").Value = "Región de Arica y Parinacota"
$ws.Range("P518").Value = 361
$ws.Range("Q518").Value = 40
$ws.Range("R518").Value = "Hortaliza"
